# Updated symbol list on Fri Dec 23 10:57:37 UTC 2022 with GitHub Actions
# Apply the latest crypto price/volume refresh to the "cryptos" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as text, not numbers, so re-assign them
# with a leading apostrophe to force text entry and preserve exact formatting
# (trailing zeros, significant digits, etc.) just like the original data.
$ws.Range("D2").Value  = "'246.21"
$ws.Range("D3").Value  = "'22.00"
$ws.Range("D4").Value  = "'5.430"
$ws.Range("D6").Value  = "'3.386"
$ws.Range("D7").Value  = "'6.351"
$ws.Range("D8").Value  = "'0.8087"
$ws.Range("D9").Value  = "'0.9690"
$ws.Range("D10").Value = "'0.1427"
$ws.Range("D11").Value = "'0.07480"
$ws.Range("D12").Value = "'0.03263"
$ws.Range("D13").Value = "'0.02999"
$ws.Range("D14").Value = "'4.141"
$ws.Range("D15").Value = "'0.09404"
$ws.Range("D16").Value = "'0.001589"
$ws.Range("D17").Value = "'0.04821"
$ws.Range("D18").Value = "'0.0005887"
$ws.Range("D19").Value = "'0.006097"
$ws.Range("D20").Value = "'0.004113"
$ws.Range("D21").Value = "'0.0009975"
$ws.Range("D24").Value = "'2.232"
$ws.Range("D27").Value = "'0.0003398"
$ws.Range("D40").Value = "'0.03876"
$ws.Range("D41").Value = "'0.006689"

$ws.Range("D43").Value = "'0.002550"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.006396"
$ws.Range("D45").Value = "'0.00005609"

$ws.Range("D47").Value = "'0.4198"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"

$ws.Range("D48").Value = "'0.1464"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D50").Value = "'0.01009"
